# Atualização de bases das ligas, do dia: 28-05-2024 às 19:13
#
# The refreshed odds feed re-ordered the match rows for this matchday
# (298-306). Column A (the running index 296-304) stays put; the rest of
# each row's data (B:AD - match id, teams, result, odds, ...) moves to a
# different row per the mapping below. We snapshot every source row's
# values first (so overlapping/cyclic moves don't clobber data we still
# need to read), then write the snapshots into their destination rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First column (B) .. last column (AD) = columns 2..30
$firstCol = 2
$lastCol = 30

# destination row -> source row (i.e. "this row's new content used to live
# in that row")
$rowMap = @{
    298 = 301
    299 = 303
    300 = 306
    301 = 305
    302 = 304
    303 = 298
    304 = 302
    305 = 300
    306 = 299
}

# Snapshot B:AD for every row involved (sources == destinations == 298..306)
# before writing anything, so cyclic permutations don't read already
# overwritten cells.
$snapshots = @{}
foreach ($srcRow in $rowMap.Values | Sort-Object -Unique) {
    $vals = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $vals[$c] = $ws.Cells.Item($srcRow, $c).Value2
    }
    $snapshots[$srcRow] = $vals
}

# Now write each destination row from its snapshot.
foreach ($destRow in $rowMap.Keys | Sort-Object) {
    $srcRow = $rowMap[$destRow]
    $vals = $snapshots[$srcRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $vals[$c]
    }
}
